$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (F column) values re-pulled / recalculated from source data.
$updates = @{
    2  = -4
    3  = -1
    4  = 0
    6  = 4
    7  = 4
    11 = 6
    13 = 2
    14 = 3
    15 = -2
    16 = -1
    17 = 7
    18 = 2
    19 = 11
    20 = -2
    23 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
